$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the annotation count for "huddle/gather" (row 8) from column G
# (similar_but_not_synonym) to column C (semrel_antonym)
$ws.Range("G8").ClearContents()
$ws.Range("C8").Value = 1

# Move the annotation count for "huddle/football" (row 9) from column G
# (similar_but_not_synonym) to column J (other)
$ws.Range("G9").ClearContents()
$ws.Range("J9").Value = 1

# Update the active selection to reflect the latest cursor position
$ws.Range("G21").Select()
